# Swap the "lat" (C) and "lng" (D) column values for the named-region rows
# (rows 11-166) on the "geolocation" sheet, where the coordinates had been
# stored in the wrong columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("geolocation")

$firstRow = 11
$lastRow = 166

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)   # Column C (lat)
    $dCell = $ws.Cells.Item($r, 4)   # Column D (lng)

    $cVal = $cCell.Value()
    $dVal = $dCell.Value()

    $cCell.Value = $dVal
    $dCell.Value = $cVal
}
